$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (default/no explicit formatting) used to restore cell style
# after forcing a text NumberFormat, so numeric-looking strings (e.g. "1.003")
# are stored as text instead of being parsed into numbers.
$plainStyle = $ws.Range("B2").Style

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $plainStyle
}

Set-TextValue "D2" '26.134.48'
$ws.Range("E2").Value = '  -0.56%  '
Set-TextValue "D3" '1.667.04'
$ws.Range("E3").Value = '  -1.18%  '
$ws.Range("E4").Value = '  -0.46%  '
Set-TextValue "D5" '209.65'
$ws.Range("E5").Value = '  -3.64%  '
Set-TextValue "D6" '0.5208'
$ws.Range("E6").Value = '  -2.16%  '
$ws.Range("E7").Value = '  -0.42%  '
Set-TextValue "D8" '0.2616'
$ws.Range("E8").Value = '  -3.56%  '
Set-TextValue "D9" '0.06323'
$ws.Range("E9").Value = '  -1.36%  '
Set-TextValue "D10" '21.10'
$ws.Range("E10").Value = '  -2.69%  '
Set-TextValue "D11" '0.07524'
$ws.Range("E11").Value = '  -2.17%  '
Set-TextValue "D12" '1.674.88'
$ws.Range("E12").Value = '  -0.49%  '
Set-TextValue "D13" '4.428'
$ws.Range("E13").Value = '  -2.25%  '
Set-TextValue "D14" '0.5495'
Set-TextValue "D15" '66.34'
$ws.Range("E15").Value = '  -0.90%  '
Set-TextValue "D16" '0.000007940'
$ws.Range("E16").Value = '  -5.14%  '
Set-TextValue "D17" '26.138.97'
$ws.Range("E17").Value = '  -0.68%  '
$ws.Range("E18").Value = '  -0.48%  '
Set-TextValue "D19" '4.720'
$ws.Range("E19").Value = '  -3.56%  '
Set-TextValue "D20" '186.39'
$ws.Range("E20").Value = '  -3.12%  '
$ws.Range("E21").Value = '  -5.28%  '
Set-TextValue "D22" '6.177'
$ws.Range("E22").Value = '  -1.33%  '
$ws.Range("E23").Value = '  -0.39%  '
Set-TextValue "D24" '149.27'
$ws.Range("E24").Value = '  +0.02%  '
Set-TextValue "D25" '0.1246'
$ws.Range("E25").Value = '  -2.07%  '
Set-TextValue "D26" '7.488'
$ws.Range("E26").Value = '  -4.46%  '
Set-TextValue "D27" '15.83'
$ws.Range("E27").Value = '  +0.03%  '
Set-TextValue "D28" '0.06366'
$ws.Range("E28").Value = '  +2.25%  '
$ws.Range("E29").Value = '  -1.95%  '
Set-TextValue "D30" '1.273'
$ws.Range("E30").Value = '  -3.85%  '
$ws.Range("E31").Value = '  -2.89%  '
Set-TextValue "D32" '3.410'
Set-TextValue "D33" '1.641'
$ws.Range("E33").Value = '  -2.91%  '
Set-TextValue "D34" '1.003'
$ws.Range("E34").Value = '  -2.72%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue "D35" '0.6021'
$ws.Range("E35").Value = '  -2.55%  '
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue "D36" '2.408'
$ws.Range("E36").Value = '  -0.90%  '
Set-TextValue "D37" '2.746'
$ws.Range("E37").Value = '  -0.37%  '
Set-TextValue "D38" '1.108.93'
$ws.Range("E38").Value = '  -0.03%  '
Set-TextValue "D39" '6.109'
$ws.Range("E39").Value = '  -2.05%  '
Set-TextValue "D40" '0.01616'
$ws.Range("E40").Value = '  -1.33%  '
Set-TextValue "D41" '0.8688'
$ws.Range("E41").Value = '  -3.18%  '
$ws.Range("E42").Value = '  -0.85%  '
Set-TextValue "D43" '100.04'
$ws.Range("E43").Value = '  -0.84%  '
Set-TextValue "D44" '1.819.50'
$ws.Range("E44").Value = '  -1.06%  '
Set-TextValue "D45" '0.00000000108'
$ws.Range("E45").Value = '  -5.29%  '
$ws.Range("E46").Value = '  -4.13%  '
$ws.Range("E47").Value = '  -0.73%  '
Set-TextValue "D48" '8.029'
$ws.Range("E48").Value = '  -0.91%  '
Set-TextValue "D49" '0.05231'
$ws.Range("E49").Value = '  -0.95%  '
Set-TextValue "D50" '0.4245'
$ws.Range("E50").Value = '  -1.11%  '
Set-TextValue "D51" '5.923'
$ws.Range("E51").Value = '  -2.33%  '
